$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 141 (pushes existing rows 141-156 down to 142-157)
$ws.Rows(141).Insert()

# Populate the newly inserted row 141 with the new weekly price record
$ws.Cells.Item(141, 1).Value = 7
$ws.Cells.Item(141, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(141, 3).Value = 'Ñuble'
$ws.Cells.Item(141, 4).Value = 45267
$ws.Cells.Item(141, 5).Value = 16
$ws.Cells.Item(141, 6).Value = 100112044
$ws.Cells.Item(141, 7).Value = 'Perejil'
$ws.Cells.Item(141, 8).Value = 'Sin especificar'
$ws.Cells.Item(141, 9).Value = 'Primera'
$ws.Cells.Item(141, 10).Value = 150
$ws.Cells.Item(141, 11).Value = 1500
$ws.Cells.Item(141, 12).Value = 1500
$ws.Cells.Item(141, 13).Value = 1500
$ws.Cells.Item(141, 14).Value = '$/atado 0,5 a 1 kilo'
$ws.Cells.Item(141, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(141, 16).Value = 1500
$ws.Cells.Item(141, 17).Value = 1
$ws.Cells.Item(141, 18).Value = 'Hortaliza'
